# Added DRA and ENW test cases for sanity job.
# This script adds a new ENW045 test case row (row 42) to the "Test Cases"
# worksheet, matching the layout/formatting of the row above it (row 41),
# and updates the view/selection state accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (and, transiently, the values) of the last existing
# data row down into the new row so the new row visually matches the rest
# of the table (borders, fills, fonts, alignment all come along for free).
$ws.Range("A41:E41").Copy($ws.Range("A42:E42"))

# Overwrite with the actual new test-case content for row 42.
$ws.Range("A42").Value = "ENW045"
$ws.Range("B42").Value = "OPQA-2015||OPQA-3650"
$ws.Range("C42").Value = "Verify that User is able to sign-into EndNote Web using STeAM."
$ws.Range("D42").Value = "Y"

# Row 42 is a single-line row (unlike the taller, wrapped row 41 above it).
$ws.Rows.Item(42).RowHeight = 15.75

# Update the window scroll position / selection like a user would have left
# it after adding the row at the bottom of the sheet.
$win = $excel.ActiveWindow
$win.ScrollRow = 35
$win.ScrollColumn = 1
$null = $ws.Range("D43").Select()
